# Insert a new weekly price record as row 268 ("Terminal Hortofrutícola Agro
# Chillán" / "Apio"), pushing the existing rows 268-312 down to 269-313.
# Excel's Rows.Insert() shifts all cell content (and formatting) of the
# rows below down by one, exactly mirroring the OOXML diff (new <row r="313">
# appended, dimension becomes A1:R313, and each prior row N becomes N+1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

$ws.Cells.Item(268, 1).Value = 7
$ws.Cells.Item(268, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(268, 3).Value = 'Ñuble'
$ws.Cells.Item(268, 4).Value = 45015
$ws.Cells.Item(268, 5).Value = 16
$ws.Cells.Item(268, 6).Value = 100112017
$ws.Cells.Item(268, 7).Value = 'Apio'
$ws.Cells.Item(268, 8).Value = 'Americana (o)'
$ws.Cells.Item(268, 9).Value = 'Primera'
$ws.Cells.Item(268, 10).Value = 150
$ws.Cells.Item(268, 11).Value = 8000
$ws.Cells.Item(268, 12).Value = 8000
$ws.Cells.Item(268, 13).Value = 8000
$ws.Cells.Item(268, 14).Value = '$/docena de matas'
$ws.Cells.Item(268, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(268, 16).Value = 1333
$ws.Cells.Item(268, 17).Value = 6
$ws.Cells.Item(268, 18).Value = 'Hortaliza'
